# Insert 3 new weekly price rows for "Ají" (chili pepper) into the
# "Mercado Mayorista Lo Valledor de Santiago" dataset, right after the
# existing row 1204, shifting the remaining rows (old 1205:1274) down to
# (1208:1277). This matches the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 1205:1274 down by three rows.
$ws.Range("A1205:A1207").EntireRow.Insert()

# Values constant across every data row in this single-market/category file.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 100112021
$categoria = "Ají"
$clasif    = "Hortaliza"

# New rows (Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm,
#           Unidad, Origen, PrecioKg, KgOUnidades)
$newRows = @(
    @{ Row=1205; Fecha=45041; Variedad="Americana (o)";   Calidad="Primera"; Vol=400; PMin=16000; PMax=18000; PProm=16850; Unidad="`$/caja 25 kilos"; Origen="Provincia de Limarí";   PKg=674; Kg=25 },
    @{ Row=1206; Fecha=45041; Variedad="Americana (o)";   Calidad="Primera"; Vol=290; PMin=15000; PMax=16000; PProm=15586; Unidad="`$/saco 25 kilos"; Origen="Región Metropolitana";  PKg=623; Kg=25 },
    @{ Row=1207; Fecha=45041; Variedad="Cacho cabra rojo"; Calidad="Primera"; Vol=200; PMin=18000; PMax=19000; PProm=18400; Unidad="`$/saco 25 kilos"; Origen="Región Metropolitana";  PKg=736; Kg=25 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $catId
    $ws.Cells.Item($row, 7).Value  = $categoria
    $ws.Cells.Item($row, 8).Value  = $r.Variedad
    $ws.Cells.Item($row, 9).Value  = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Vol
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $r.Unidad
    $ws.Cells.Item($row, 15).Value = $r.Origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $r.Kg
    $ws.Cells.Item($row, 18).Value = $clasif
}
